$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SuperAdmin")
$ws2 = $wb.Worksheets.Item("Customer")

# --- SuperAdmin (sheet1) cell updates ---
$ws1.Range("A2").Value = "judecisla@gmail.com"
$ws1.Range("B3").Value = "Admin123"

# A4 needs the same border-only style that A3/B3 already carry (s="1").
$ws1.Range("A3").Copy() | Out-Null
$ws1.Range("A4").PasteSpecial(-4122) | Out-Null

# B4 gets a password value and loses its border while picking up an
# (empty) fill-applied style.
$ws1.Range("B4").Borders.LineStyle = -4142
$ws1.Range("B4").Interior.ColorIndex = -4142
$ws1.Range("B4").Value = "Hello@1234"

$ws1.Range("A5").Value = "judecisla@gmail.com"

# New dropdown validation for the password column, mirroring the
# existing username validation on column A.
$ws1.Range("B2").Validation.Add(3, 1, 1, "=`$B`$3:`$B`$10")

# --- Customer (sheet2) cell updates ---
$ws2.Range("A2").Value = "bartryfry@macr2.com"

# --- Selection / active tab bookkeeping ---
$ws2.Range("A2").Select() | Out-Null
$ws1.Activate()
$ws1.Range("D3").Select() | Out-Null

Write-Output "done"
